# Auto-generated Excel COM-interop script
# Applies numeric cell value updates to "Output_flows" and "Input_flows" sheets
# as described by the source diff (updated mass-flow values for new input files).

$wb = $excel.ActiveWorkbook

$wsOutput = $wb.Worksheets.Item("Output_flows")
$wsInput  = $wb.Worksheets.Item("Input_flows")

# --- Output_flows sheet updates ---
$wsOutput.Range("C2").Value = [double]"1.794957884292373E-16"
$wsOutput.Range("E2").Value = [double]"6.434694712315917E-13"
$wsOutput.Range("G2").Value = [double]"4.31596139712192E-14"
$wsOutput.Range("I2").Value = [double]"3.720554016690005E-14"
$wsOutput.Range("K2").Value = [double]"1.442425064066983E-17"
$wsOutput.Range("C7").Value = [double]"1.902341373096434E-08"
$wsOutput.Range("F7").Value = [double]"6.905689668514041E-13"
$wsOutput.Range("G7").Value = [double]"4.574164108406879E-05"
$wsOutput.Range("I7").Value = [double]"3.943136437198241E-05"
$wsOutput.Range("J7").Value = [double]"2.436143591642356E-07"
$wsOutput.Range("C12").Value = [double]"1.776212128760545E-14"
$wsOutput.Range("E12").Value = [double]"1.411839116737883E-07"
$wsOutput.Range("I12").Value = [double]"1.840849088433856E-11"
$wsOutput.Range("J12").Value = [double]"2.464427933841764E-11"
$wsOutput.Range("C13").Value = [double]"2.197882565141099E-15"
$wsOutput.Range("D13").Value = [double]"1.628873772526735E-15"
$wsOutput.Range("E13").Value = [double]"6.515521138168181E-08"
$wsOutput.Range("I13").Value = [double]"2.277864254506563E-12"
$wsOutput.Range("J13").Value = [double]"5.964223951167562E-12"
$wsOutput.Range("C14").Value = [double]"3.057452297684681E-16"
$wsOutput.Range("D14").Value = [double]"9.063639591008744E-16"
$wsOutput.Range("E14").Value = [double]"7.773908239750561E-07"
$wsOutput.Range("I14").Value = [double]"3.168714020126818E-13"
$wsOutput.Range("J14").Value = [double]"4.812670126678333E-12"
$wsOutput.Range("C15").Value = [double]"1.275745093320558E-35"
$wsOutput.Range("D15").Value = [double]"3.214591338978486E-34"
$wsOutput.Range("E15").Value = [double]"2.958428400363186E-23"
$wsOutput.Range("I15").Value = [double]"1.322169888430996E-32"
$wsOutput.Range("K15").Value = [double]"2.616586762609753E-30"
$wsOutput.Range("C17").Value = [double]"4.436928844264685E-06"
$wsOutput.Range("F17").Value = [double]"1.411112918236624E-07"
$wsOutput.Range("I17").Value = [double]"0.01839356073783827"
$wsOutput.Range("J17").Value = [double]"0.02463056265127432"
$wsOutput.Range("C18").Value = [double]"5.493360424988749E-07"
$wsOutput.Range("D18").Value = [double]"1.628494866627871E-08"
$wsOutput.Range("F18").Value = [double]"6.514637068389206E-08"
$wsOutput.Range("I18").Value = [double]"0.002277306266979654"
$wsOutput.Range("J18").Value = [double]"0.005963309646617993"
$wsOutput.Range("C19").Value = [double]"7.643605064409196E-08"
$wsOutput.Range("D19").Value = [double]"9.063610888877153E-09"
$wsOutput.Range("F19").Value = [double]"7.773916912968154E-07"
$wsOutput.Range("I19").Value = [double]"0.0003168703374407138"
$wsOutput.Range("J19").Value = [double]"0.004812674607625756"
$wsOutput.Range("C20").Value = [double]"3.189362857265561E-27"
$wsOutput.Range("D20").Value = [double]"3.214591464932968E-27"
$wsOutput.Range("F20").Value = [double]"2.958428400395949E-23"
$wsOutput.Range("I20").Value = [double]"1.322169939821102E-23"
$wsOutput.Range("K20").Value = [double]"2.616586757563879E-21"

# --- Input_flows sheet updates ---
$wsInput.Range("C2").Value = [double]"3.345957855737681E-14"
$wsInput.Range("C7").Value = [double]"4.860455215802391E-05"
$wsInput.Range("C12").Value = [double]"1.387926862542485E-11"
$wsInput.Range("C13").Value = [double]"1.52679144208826E-12"
$wsInput.Range("C14").Value = [double]"1.785125767952804E-13"
$wsInput.Range("C15").Value = [double]"1.716622137830005E-36"
$wsInput.Range("C17").Value = [double]"0.0135678695167558"
$wsInput.Range("C18").Value = [double]"0.001515150177511291"
$wsInput.Range("C19").Value = [double]"0.0001785037723130402"
$wsInput.Range("C20").Value = [double]"1.716622204665843E-27"
$wsInput.Range("D27").Value = [double]"3.674634116665635E-05"
$wsInput.Range("D32").Value = [double]"9.609770644213841E-11"
$wsInput.Range("D33").Value = [double]"1.420554360296129E-11"
$wsInput.Range("D34").Value = [double]"3.198658393985607E-12"
$wsInput.Range("D37").Value = [double]"0.0175204350142379"
$wsInput.Range("D38").Value = [double]"0.003846332529388074"
$wsInput.Range("D39").Value = [double]"0.002639304320130874"
$wsInput.Range("C42").Value = [double]"3.832048609676356E-32"
$wsInput.Range("C47").Value = [double]"3.648163975323471E-23"
$wsInput.Range("C52").Value = [double]"2.061359206047153E-29"
$wsInput.Range("C53").Value = [double]"2.383607060060894E-30"
$wsInput.Range("C54").Value = [double]"2.589049278097644E-31"
$wsInput.Range("C55").Value = [double]"2.629813324448865E-30"
$wsInput.Range("C57").Value = [double]"2.060897853042382E-20"
$wsInput.Range("C58").Value = [double]"2.383379746319185E-21"
$wsInput.Range("C59").Value = [double]"2.589037495295418E-22"
$wsInput.Range("C60").Value = [double]"2.629813144294536E-21"
$wsInput.Range("F67").Value = [double]"8.474995136565549E-08"
$wsInput.Range("F72").Value = [double]"5.668165732666766E-12"
$wsInput.Range("F73").Value = [double]"1.353824524739553E-12"
$wsInput.Range("F74").Value = [double]"8.862609077492111E-13"
$wsInput.Range("F77").Value = [double]"0.01189449325650512"
$wsInput.Range("F78").Value = [double]"0.002879694287043307"
$wsInput.Range("F79").Value = [double]"0.002311822353151412"
